# Updated cryptos list on Thu Jun  6 22:42:21 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values for the
# cryptocurrency rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.780.60"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "3.805.17"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "708.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("D7").Value = "3.805.64"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "4.449.86"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "3.796.54"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "70.826.59"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "494.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.75%  "
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("D28").Value = "3.958.44"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.173"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("D37").Value = "3.777.77"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("E45").Value = "  +5.60%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "165.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "425.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.295"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.86%  "
